$d = $word.ActiveDocument

# "What we like" bullet list updates
$d.Content.Find.Execute("Well-designed Halloween theme", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Exciting and well-designed Halloween theme", 2)

$d.Content.Find.Execute("Innovative gameplay elements", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Many bonus features", 2)

$d.Content.Find.Execute("High roller slot game", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Innovative gameplay elements", 2)

# "What we don't like" bullet list updates
$d.Content.Find.Execute("Slightly lower RTP than average", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Slightly lower than average RTP", 2)

$d.Content.Find.Execute("Not suitable for low-risk, low-reward players", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "High minimum and maximum bets", 2)

# Meta/summary italic paragraph rewrite
$d.Content.Find.Execute("Read our review of Jack O’Lantern vs The Headless Horseman, a Halloween-themed slot game with a huge jackpot of €4,000,000 and many bonus features. Play for free.", `
                         $true, $true, $false, $false, $false, `
                         $true, 1, $false, "Read a review of Jack O’Lantern vs The Headless Horseman slot game and play for free. Enjoy a thrilling Halloween-themed experience with a huge €4,000,000 jackpot.", 2)
